$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 123.26667
$ws.Range("I33").Value = 164.28572
$ws.Range("J33").Value = 87.375
$ws.Range("K33").Value = 164.28572
$ws.Range("L33").Value = 87.375
$ws.Range("M33").Value = 64.71428
$ws.Range("N33").Value = -545.375
$ws.Range("H69").Value = 4000
$ws.Range("I69").Value = 4000
$ws.Range("K69").Value = 12000
$ws.Range("M69").Value = -11126
$ws.Range("H72").Value = 4000
$ws.Range("I72").Value = 4000
$ws.Range("K72").Value = 36000
$ws.Range("M72").Value = -31632
$ws.Range("H86").Value = 7424.8335
$ws.Range("I86").Value = 6250
$ws.Range("K86").Value = 6250
$ws.Range("M86").Value = -5127
$ws.Range("H89").Value = 7424.8335
$ws.Range("I89").Value = 6250
$ws.Range("K89").Value = 31250
$ws.Range("M89").Value = -25634
$ws.Range("H98").Value = 1685.2941
$ws.Range("J98").Value = 1629.5
$ws.Range("L98").Value = 1629.5
$ws.Range("N98").Value = -4625.5
$ws.Range("H106").Value = 3266.1428
$ws.Range("I106").Value = 3266.1428
$ws.Range("K106").Value = 3266.1428
$ws.Range("M106").Value = -2635.1428
$ws.Range("H113").Value = 6549.4
$ws.Range("I113").Value = 5125
$ws.Range("K113").Value = 5125
$ws.Range("M113").Value = -1871
$ws.Range("H122").Value = 1685.2941
$ws.Range("J122").Value = 1629.5
$ws.Range("L122").Value = 4888.5
$ws.Range("N122").Value = -9788.5
$ws.Range("H125").Value = 500499.5
$ws.Range("J125").Value = 999999
$ws.Range("L125").Value = 8999991
$ws.Range("N125").Value = -9004911
$ws.Range("H130").Value = 74000
$ws.Range("J130").Value = 74000
$ws.Range("L130").Value = 74000
$ws.Range("N130").Value = -84040

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2782.92
$ws.Range("I32").Value = 1233.95
$ws.Range("K32").Value = 1233.95
$ws.Range("M32").Value = -946.95
$ws.Range("H101").Value = 48749.75
$ws.Range("J101").Value = 48749.75
$ws.Range("L101").Value = 48749.75
$ws.Range("N101").Value = -55239.75
$ws.Range("H131").Value = 60999.2
$ws.Range("J131").Value = 60999.2
$ws.Range("L131").Value = 60999.2
$ws.Range("N131").Value = -71079.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 3551.8
$ws.Range("I20").Value = 3229.5
$ws.Range("J20").Value = 3766.6667
$ws.Range("K20").Value = 3229.5
$ws.Range("L20").Value = 3766.6667
$ws.Range("M20").Value = -2982.5
$ws.Range("N20").Value = -4260.6667
$ws.Range("H81").Value = 42717.332
$ws.Range("J81").Value = 42717.332
$ws.Range("L81").Value = 42717.332
$ws.Range("N81").Value = -44839.332
$ws.Range("H84").Value = 42717.332
$ws.Range("J84").Value = 42717.332
$ws.Range("L84").Value = 128151.996
$ws.Range("N84").Value = -138759.996

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1386
$ws.Range("I16").Value = 1608.5
$ws.Range("J16").Value = 496
$ws.Range("K16").Value = 1608.5
$ws.Range("L16").Value = 496
$ws.Range("M16").Value = -1321.5
$ws.Range("N16").Value = -1070
$ws.Range("H62").Value = 3750
$ws.Range("J62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("H65").Value = 3750
$ws.Range("J65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H105").Value = 1844.75
$ws.Range("I105").Value = 1844.75
$ws.Range("K105").Value = 1844.75
$ws.Range("M105").Value = -97.75
$ws.Range("H107").Value = 355.46155
$ws.Range("I107").Value = 236.55556
$ws.Range("K107").Value = 236.55556
$ws.Range("M107").Value = 1683.44444
$ws.Range("H113").Value = 1386
$ws.Range("I113").Value = 1608.5
$ws.Range("J113").Value = 496
$ws.Range("K113").Value = 1608.5
$ws.Range("L113").Value = 496
$ws.Range("M113").Value = 561.5
$ws.Range("N113").Value = -4836
$ws.Range("H134").Value = 3610.5938
$ws.Range("I134").Value = 2491.4285
$ws.Range("J134").Value = 4481.0557
$ws.Range("K134").Value = 7474.2855
$ws.Range("L134").Value = 13443.1671
$ws.Range("M134").Value = -4939.2855
$ws.Range("N134").Value = -18513.1671

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 1999
$ws.Range("J93").Value = 1999
$ws.Range("L93").Value = 5997
$ws.Range("N93").Value = -9741
$ws.Range("H113").Value = 1028.5
$ws.Range("J113").Value = 1358
$ws.Range("L113").Value = 4074
$ws.Range("N113").Value = -8414
$ws.Range("H130").Value = 2143
$ws.Range("I130").Value = 1029
$ws.Range("J130").Value = 2700
$ws.Range("K130").Value = 3087
$ws.Range("L130").Value = 8100
$ws.Range("M130").Value = 1933
$ws.Range("N130").Value = -18140

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 21999.285
$ws.Range("J33").Value = 21999.285
$ws.Range("L33").Value = 21999.285
$ws.Range("N33").Value = -22503.285
$ws.Range("H62").Value = 20000
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 20000
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("N92").Value = 0
$ws.Range("H97").Value = 2426.25
$ws.Range("I97").Value = 1835
$ws.Range("J97").Value = 4200
$ws.Range("K97").Value = 1835
$ws.Range("L97").Value = 4200
$ws.Range("M97").Value = -1339
$ws.Range("N97").Value = -5192
$ws.Range("H98").Value = 8179.8
$ws.Range("J98").Value = 8179.8
$ws.Range("L98").Value = 8179.8
$ws.Range("N98").Value = -14169.8
$ws.Range("H99").Value = 23911.666
$ws.Range("I99").Value = 23911.666
$ws.Range("J99").Value = 0
$ws.Range("K99").Value = 23911.666
$ws.Range("L99").Value = 0
$ws.Range("N99").Value = -21665.666
$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("N100").Value = 0
$ws.Range("H102").Value = 13737.546
$ws.Range("J102").Value = 23783
$ws.Range("L102").Value = 23783
$ws.Range("N102").Value = -27027
$ws.Range("H105").Value = 39888.332
$ws.Range("J105").Value = 151497.5
$ws.Range("L105").Value = 151497.5
$ws.Range("N105").Value = -158485.5
$ws.Range("H107").Value = 463.33334
$ws.Range("I107").Value = 463.33334
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 463.33334
$ws.Range("L107").Value = 0
$ws.Range("N107").Value = 1456.66666
$ws.Range("H122").Value = 2408.8823
$ws.Range("I122").Value = 1831.8
$ws.Range("J122").Value = 3233.2856
$ws.Range("K122").Value = 5495.4
$ws.Range("L122").Value = 9699.856800000001
$ws.Range("M122").Value = -3045.4
$ws.Range("N122").Value = -14599.8568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H50").Value = 40000
$ws.Range("J50").Value = 40000
$ws.Range("L50").Value = 40000
$ws.Range("N50").Value = -41274
$ws.Range("H55").Value = 280.45456
$ws.Range("I55").Value = 304.33334
$ws.Range("K55").Value = 304.33334
$ws.Range("M55").Value = -131.33334
$ws.Range("H56").Value = 42000
$ws.Range("J56").Value = 42000
$ws.Range("L56").Value = 42000
$ws.Range("N56").Value = -43382
$ws.Range("H63").Value = 32271.25
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 32271.25
$ws.Range("K63").Value = 0
$ws.Range("M63").Value = 32271.25
$ws.Range("N63").Value = -33769.25
$ws.Range("H66").Value = 32271.25
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 32271.25
$ws.Range("K66").Value = 0
$ws.Range("M66").Value = 96813.75
$ws.Range("N66").Value = -104301.75
$ws.Range("H106").Value = 0
$ws.Range("J106").Value = 0
$ws.Range("N106").Value = 0

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H51").Value = 54961.668
$ws.Range("I51").Value = 54770
$ws.Range("K51").Value = 54770
$ws.Range("M51").Value = -54260
$ws.Range("H52").Value = 7300.75
$ws.Range("I52").Value = 7300.75
$ws.Range("J52").Value = 0
$ws.Range("K52").Value = 7300.75
$ws.Range("L52").Value = 0
$ws.Range("N52").Value = -7074.75
$ws.Range("H81").Value = 3356.8572
$ws.Range("J81").Value = 3499
$ws.Range("L81").Value = 6998
$ws.Range("N81").Value = -9120
$ws.Range("H84").Value = 3356.8572
$ws.Range("J84").Value = 3499
$ws.Range("L84").Value = 34990
$ws.Range("N84").Value = -45598
$ws.Range("H113").Value = 639.8
$ws.Range("J113").Value = 500
$ws.Range("L113").Value = 1500
$ws.Range("N113").Value = -5840
